$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# Correct the erroneous values in row 10 (2018)
$ws.Range("C10").Value = 28103000000
$ws.Range("D10").Value = 2974000000
$ws.Range("E10").Value = 2163000000

# Update the active selection on the sheet
$ws.Range("H7").Select()
